$wb = $excel.ActiveWorkbook

# --- "Backlog" sheet: reword the bauble/track-section linking story, and
# mark the two rows (118, 120) that now roll up into the completion count.
$bl = $wb.Worksheets.Item("Backlog")

$bl.Range("B120").Value = "Once the file is completely loaded and deserialized, the track tool must go through the track sections and link to bauble objects based on IDs."
$bl.Range("C118").Value = "X"
$bl.Range("C120").Value = "X"

# --- "Stories 6" sheet: add the new "Finalising a track section..." story
# after the existing rows, and the new "Upgrade to Unity 5" story before
# the existing rows.
$s6 = $wb.Worksheets.Item("Stories 6")

$s6.Range("A3").Value = "Finalising a track section should not be so slow if it doesn't need to adjust the terrain."
$s6.Rows.Item(3).RowHeight = 39

$s6.Rows.Item(1).Insert()
$s6.Range("A1").Value = "Upgrade to Unity 5"
$s6.Rows.Item(1).RowHeight = 39

# --- view/selection bookkeeping to match the edited workbook's cursor
# state (Stories 6 selection first, then re-activate Backlog so it keeps
# being the selected tab).
$s6.Range("A2").Select()

$bl.Activate()
$bl.Range("C121").Select()
